$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" -----------
# "Overview" sheet tracks status per-language in columns E (zh-cn) and F (de-de).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# "zh-cn" and "de-de" sheets track status in column C.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "In Translation"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "In Translation"

# --- Re-fit the status columns now that the text is shorter -----------------
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZh.Columns.Item(3).ColumnWidth = $newWidth
$wsDe.Columns.Item(3).ColumnWidth = $newWidth
